# repull data, push all data, mean calculation
# Update column F (dSF) values for specific rows to reflect the recalculated
# "final" delta-S figures (previously these mirrored column E / dS0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -6
    4  = -4
    16 = 0
    18 = 1
    19 = -1
    20 = -2
    24 = 1
    25 = 0
    27 = 3
    33 = -3
    36 = -1
    39 = 0
    42 = 1
    43 = 3
    47 = -4
    50 = -6
    51 = -7
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
